$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns retain text formatting
# (values like "10.00" or "0.0000180" must not be auto-converted to numbers)
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.553.25"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "3.507.37"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D5").Value = "586.48"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "133.04"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "3.506.42"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").Value = "7.15"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "4.102.65"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "27.86"
$ws.Range("E14").Value = "  +4.01%  "
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "3.503.78"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "64.562.52"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "14.26"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("D22").Value = "393.36"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "0.579"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "3.646.36"
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("D25").Value = "74.19"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "0.0000110"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "7.48"
$ws.Range("E29").Value = "  -3.47%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").Value = "3.510.71"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "23.99"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").Value = "1.61"
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("D38").Value = "5.28"
$ws.Range("E38").Value = "  +5.34%  "
$ws.Range("D39").Value = "171.22"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "6.99"
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("D41").Value = "0.0816"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "0.813"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").Value = "26.52"
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "42.17"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D46").Value = "1.21"
$ws.Range("E46").Value = "  -2.93%  "
$ws.Range("D47").Value = "4.41"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").Value = "2.467.00"
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("D50").Value = "6.90"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "0.909"
$ws.Range("E51").Value = "  +5.42%  "

# Restore default (Normal) style so no stray number-format styling is introduced
$priceVolumeRange.Style = "Normal"
